# Penalty Reward System (unfinished) — shift the forecast week/date rows
# forward by one week and zero-out most of the MyForecast values, then
# patch a handful of stats on the Summary sheet to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet: columns B (Week_Start_Date) and D (MyForecast) ---
# Force text format on column B before assigning so the date-shaped strings are
# kept as literal text instead of being auto-converted to Excel date serials.

$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2025-01-12"
$ws1.Range("D2").Value = 1

$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "2025-01-19"
$ws1.Range("D3").Value = 1

$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2025-01-26"
$ws1.Range("D4").Value = 0

$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2025-02-02"
$ws1.Range("D5").Value = 0

$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2025-02-09"
$ws1.Range("D6").Value = 0

$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2025-02-16"
$ws1.Range("D7").Value = 0

$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2025-02-23"
$ws1.Range("D8").Value = 0

$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "2025-03-02"
$ws1.Range("D9").Value = 0

$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("D10").Value = 0

$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("D11").Value = 0

$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("D12").Value = 0

$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("D13").Value = 0

$ws1.Range("B14").NumberFormat = "@"
$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("D14").Value = 0

$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("D15").Value = 0

$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("D16").Value = 0

$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "2025-04-27"
$ws1.Range("D17").Value = 0

# --- "Summary" sheet: a handful of recomputed/edited metric values ---

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"

$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "84"

$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "18"

$ws2.Range("B8").Value = "1884 units"

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "7"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "4"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "2"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "1"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "0"

$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-02-02"
